$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the three question labels in column A (drop the " is a problem" suffix)
$ws.Range("A2").Value = "Income inequality in [Country]"
$ws.Range("A3").Value = "Climate change"
$ws.Range("A4").Value = "Global poverty"

# Overwrite figures in column B with the corrected (higher-precision) past values
$ws.Range("B2").Value = 0.65434195987189
$ws.Range("B3").Value = 0.684162620459082
$ws.Range("B4").Value = 0.657026216062907
